$wb = $excel.ActiveWorkbook

# --- "About" sheet: refresh the "last updated" date shown in C1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = "3/28/2024"

# Scroll the About sheet so row 6 is at the top of the view (best effort;
# mirrors the saved scroll position from the authored workbook).
$wsAbout.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1

# --- "FPIEBP" sheet: re-prioritize hard coal's production/imports/exports ---
$wsFPIEBP = $wb.Worksheets.Item("FPIEBP")
$wsFPIEBP.Activate()

$wsFPIEBP.Range("B3").Value = 1
$wsFPIEBP.Range("C3").Value = 3
$wsFPIEBP.Range("D3").Value = 2

# Move the active selection to E3, matching the author's last cursor position.
[void]$wsFPIEBP.Range("E3").Select()
